$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.047.98'
$ws.Range("E2").Value = '  +2.08%  '

$ws.Range("D3").Value = '1.673.66'
$ws.Range("E3").Value = '  +2.72%  '

$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.10'
$ws.Range("E5").Value = '  +1.41%  '

$ws.Range("E6").Value = '  +1.83%  '

$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.251'
$ws.Range("E8").Value = '  +1.91%  '

$ws.Range("E9").Value = '  +1.25%  '

$ws.Range("E10").Value = '  +4.98%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0889'
$ws.Range("E11").Value = '  +4.48%  '

$ws.Range("D12").Value = '1.909.58'
$ws.Range("E12").Value = '  +2.64%  '

$ws.Range("D13").Value = '1.674.31'
$ws.Range("E13").Value = '  +2.85%  '

$ws.Range("E14").Value = '  +1.26%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '65.81'
$ws.Range("E15").Value = '  +2.73%  '

$ws.Range("E16").Value = '  +1.74%  '

$ws.Range("D17").Value = '27.067.42'
$ws.Range("E17").Value = '  +2.07%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '234.97'
$ws.Range("E18").Value = '  -0.64%  '

$ws.Range("E19").Value = '  +1.55%  '

$ws.Range("E20").Value = '  -1.16%  '

$ws.Range("E21").Value = '  +0.14%  '

$ws.Range("E22").Value = '  +3.01%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.24'
$ws.Range("E23").Value = '  +2.05%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.26'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.32'
$ws.Range("E25").Value = '  -1.18%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.17'
$ws.Range("E26").Value = '  +1.05%  '

$ws.Range("E27").Value = '  +0.47%  '

$ws.Range("E28").Value = '  +1.34%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.01%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0499'
$ws.Range("E30").Value = '  +0.14%  '

$ws.Range("E31").Value = '  +1.72%  '

$ws.Range("E32").Value = '  +2.07%  '

$ws.Range("D33").Value = '1.456.61'
$ws.Range("E33").Value = '  -4.18%  '

$ws.Range("E34").Value = '  +5.24%  '

$ws.Range("E35").Value = '  +5.41%  '

$ws.Range("E36").Value = '  -0.37%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.571'
$ws.Range("E37").Value = '  +0.20%  '

$ws.Range("E38").Value = '  +6.92%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0170'
$ws.Range("E39").Value = '  +1.46%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.08'
$ws.Range("E40").Value = '  +3.39%  '

$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.10%  '

$ws.Range("B42").Value = 'WEMIXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.01'
$ws.Range("E42").Value = '  +11.85%  '

$ws.Range("E43").Value = '  +2.93%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '65.91'
$ws.Range("E44").Value = '  +4.54%  '

$ws.Range("D45").Value = '1.819.27'
$ws.Range("E45").Value = '  +2.73%  '

$ws.Range("E46").Value = '  +2.50%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.27'
$ws.Range("E47").Value = '  -0.09%  '

$ws.Range("E48").Value = '  +1.26%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.100'
$ws.Range("E49").Value = '  +4.09%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0509'
$ws.Range("E50").Value = '  +1.55%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.65'
$ws.Range("E51").Value = '  +1.77%  '
